$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.518.96"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.911.17"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.73"
$ws.Range("E5").Value = "  -3.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.76"
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.910.53"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.501"
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.89"
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.150"
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.433"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000239"
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.91"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.396.25"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.495.15"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.64"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.915.03"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "430.23"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.16"
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.663"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.90"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.87"
$ws.Range("E24").Value = "  -2.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.04"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.98"
$ws.Range("E27").Value = "  -1.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.03"
$ws.Range("E28").Value = "  -1.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000110"
$ws.Range("E29").Value = "  +4.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.08"
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.49"
$ws.Range("E31").Value = "  -2.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.03"
$ws.Range("E32").Value = "  -3.16%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.106"
$ws.Range("E34").Value = "  -1.83%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.81"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.954"
$ws.Range("E36").Value = "  -2.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.41"
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.96"
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "48.91"
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.90"
$ws.Range("E40").Value = "  -4.90%  "
$ws.Range("B41").Value = "Arweave"
$ws.Range("C41").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "41.29"
$ws.Range("E41").Value = "  +6.84%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.114"
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.10"
$ws.Range("E43").Value = "  -2.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.268"
$ws.Range("E44").Value = "  -1.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.712.68"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0339"
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "133.60"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "353.44"
$ws.Range("E48").Value = "  +3.57%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  +16.11%  "
$ws.Range("E51").Value = "  -0.50%  "
